$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "MarketCap"
$ws.Range("G11").Select()
